# Applies the "Updated cryptos list" data refresh described by the commit diff:
# refreshed Price (D) / Volume(1h) (E) figures for every row, plus four rows
# (46-49) whose coin identity (Coin/Link in B/C) shifted position in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.970.90'
$ws.Range('E2').Value = '  +0.65%  '
# Row 3
$ws.Range('D3').Value = '1.842.55'
$ws.Range('E3').Value = '  +1.80%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  -0.03%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.67'
$ws.Range('E5').Value = '  +0.10%  '
# Row 6
$ws.Range('E6').Value = '  +2.66%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('E7').Value = '  -0.04%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.36'
$ws.Range('E8').Value = '  +5.25%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.327'
$ws.Range('E9').Value = '  +2.36%  '
# Row 10
$ws.Range('E10').Value = '  +1.65%  '
# Row 11
$ws.Range('E11').Value = '  -1.39%  '
# Row 12
$ws.Range('E12').Value = '  +1.87%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.44'
$ws.Range('E13').Value = '  +3.98%  '
# Row 14
$ws.Range('D14').Value = '1.850.42'
$ws.Range('E14').Value = '  +2.18%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.671'
$ws.Range('E15').Value = '  +0.93%  '
# Row 16
$ws.Range('E16').Value = '  +2.43%  '
# Row 17
$ws.Range('D17').Value = '35.007.90'
$ws.Range('E17').Value = '  +0.85%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.94'
$ws.Range('E18').Value = '  +0.41%  '
# Row 19
$ws.Range('E19').Value = '  +0.36%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.05'
$ws.Range('E20').Value = '  -0.07%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.16'
$ws.Range('E21').Value = '  +2.01%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.76'
$ws.Range('E22').Value = '  +2.58%  '
# Row 23
$ws.Range('E23').Value = '  -0.20%  '
# Row 24
$ws.Range('E24').Value = '  +3.45%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.23'
$ws.Range('E25').Value = '  -0.98%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.85'
$ws.Range('E26').Value = '  +1.37%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.46'
$ws.Range('E27').Value = '  +1.47%  '
# Row 28
$ws.Range('E28').Value = '  +3.36%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.73'
$ws.Range('E29').Value = '  +11.05%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.01'
$ws.Range('E30').Value = '  +0.13%  '
# Row 31
$ws.Range('E31').Value = '  +0.75%  '
# Row 32
$ws.Range('E32').Value = '  -2.16%  '
# Row 33
$ws.Range('E33').Value = '  -1.62%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.64'
$ws.Range('E34').Value = '  +22.72%  '
# Row 35
$ws.Range('E35').Value = '  +10.67%  '
# Row 36
$ws.Range('E36').Value = '  -1.60%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.756'
$ws.Range('E37').Value = '  +7.59%  '
# Row 38
$ws.Range('E38').Value = '  +9.45%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '89.98'
$ws.Range('E39').Value = '  -1.71%  '
# Row 40
$ws.Range('E40').Value = '  +3.78%  '
# Row 41
$ws.Range('D41').Value = '1.342.20'
$ws.Range('E41').Value = '  +1.74%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.50'
$ws.Range('E42').Value = '  +2.00%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.26'
$ws.Range('E43').Value = '  +1.63%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.41'
$ws.Range('E44').Value = '  -2.41%  '
# Row 45
$ws.Range('E45').Value = '  +3.77%  '
# Row 46
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0529'
$ws.Range('E46').Value = '  +3.48%  '
# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.30'
$ws.Range('E47').Value = '  +1.21%  '
# Row 48
$ws.Range('B48').Value = 'Gas'
$ws.Range('C48').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.40'
$ws.Range('E48').Value = '  +70.33%  '
# Row 49
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.022.95'
$ws.Range('E49').Value = '  +1.27%  '
# Row 50
$ws.Range('E50').Value = '  +0.16%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  -0.09%  '
